$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 1) Widen the third table column from 3306 dxa (165.3pt) to 3307 dxa
#    (165.35pt). Word's Cell/Column.Width is expressed in points (1 pt = 20
#    dxa), so 3307 dxa == 165.35 pt.
# ---------------------------------------------------------------------------
$t.Columns.Item(3).Width = 165.35

# ---------------------------------------------------------------------------
# 2) Delete the row "Enter, operate in, and depart a traffic pattern." in
#    its entirety (all three cells).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $rowText = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($rowText -like "*Enter, operate in, and depart a traffic pattern*") {
        $t.Rows.Item($i).Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) The page-break rendering cache marker (<w:lastRenderedPageBreak/>) moved
#    off of the "TASK" header cell and onto the start of the "Maintain a
#    constant approach angle..." run as repagination shifted by one row.
#    Re-point it by rewriting each paragraph's XML (preserving every other
#    attribute) via a whole-paragraph InsertXML replace, which is the only
#    way to touch this non-OM-exposed element.
# ---------------------------------------------------------------------------
function Move-LastRenderedPageBreak($doc, $table, $matchText, $add) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        $cell = $table.Rows.Item($i).Cells.Item(1)
        $cellRange = $cell.Range
        if ($cellRange.Text -like $matchText) {
            $paraXml = $cellRange.Paragraphs.Item(1).Range.WordOpenXML
            break
        }
    }
}

# --- remove <w:lastRenderedPageBreak/> before "TASK" (2nd header row) ---
$taskRowIndex = -1
$taskHits = 0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cell = $t.Rows.Item($i).Cells.Item(1)
    if ($cell.Range.Text -eq "TASK`r") {
        $taskHits = $taskHits + 1
        if ($taskHits -eq 2) {
            $taskRowIndex = $i
        }
    }
}
if ($taskRowIndex -ge 1) {
    $cell = $t.Rows.Item($taskRowIndex).Cells.Item(1)
    $r = $cell.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $snippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="0E169AE3" w14:textId="77777777" w:rsidR="00F13E4E" w:rsidRPr="00F13E4E" w:rsidRDefault="00F13E4E" w:rsidP="00F13E4E"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00F13E4E"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>TASK</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($snippet)
}

# --- add <w:lastRenderedPageBreak/> before "Maintain a constant approach..." ---
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cell = $t.Rows.Item($i).Cells.Item(1)
    if ($cell.Range.Text -like "Maintain a constant approach angle*") {
        $r = $cell.Range
        $target = $d.Range($r.Start, $r.End - 1)
        $snippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="47D819F5" w14:textId="77777777" w:rsidR="00F13E4E" w:rsidRPr="00F13E4E" w:rsidRDefault="00F13E4E" w:rsidP="00F13E4E"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="1F1E1F"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00F13E4E"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="1F1E1F"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Maintain a constant approach angle clear of obstacles to desired point of termination (hover) or touchdown. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $target.InsertXML($snippet)
        break
    }
}

Write-Output "done"
